$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the Price/Volume columns so that values
# such as "0.5900" or "1.794.66" are preserved exactly as strings and
# are not re-interpreted as numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '28.566.66'
$ws.Range('E2').Value = '  -2.29%  '
$ws.Range('D3').Value = '1.793.92'
$ws.Range('E3').Value = '  -1.87%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '231.44'
$ws.Range('E5').Value = '  -1.31%  '
$ws.Range('D6').Value = '0.5900'
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').Value = '0.2775'
$ws.Range('E8').Value = '  +1.07%  '
$ws.Range('D9').Value = '23.41'
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('D10').Value = '0.06760'
$ws.Range('E10').Value = '  -2.86%  '
$ws.Range('D11').Value = '0.07549'
$ws.Range('E11').Value = '  -1.03%  '
$ws.Range('D12').Value = '1.794.74'
$ws.Range('E12').Value = '  -2.07%  '
$ws.Range('D13').Value = '4.795'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').Value = '0.6131'
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('D15').Value = '2.036.91'
$ws.Range('E15').Value = '  -1.90%  '
$ws.Range('D16').Value = '75.68'
$ws.Range('E16').Value = '  -3.45%  '
$ws.Range('D17').Value = '0.000008891'
$ws.Range('E17').Value = '  -8.48%  '
$ws.Range('D18').Value = '28.556.87'
$ws.Range('E18').Value = '  -1.25%  '
$ws.Range('D19').Value = '5.427'
$ws.Range('E19').Value = '  -5.36%  '
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('D21').Value = '209.05'
$ws.Range('E21').Value = '  -5.67%  '
$ws.Range('D22').Value = '11.48'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('D23').Value = '6.829'
$ws.Range('E23').Value = '  -0.57%  '
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('D25').Value = '152.48'
$ws.Range('E25').Value = '  -2.23%  '
$ws.Range('D26').Value = '8.010'
$ws.Range('E26').Value = '  +1.06%  '
$ws.Range('D27').Value = '0.1265'
$ws.Range('E27').Value = '  -1.80%  '
$ws.Range('D28').Value = '16.44'
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('D29').Value = '1.418'
$ws.Range('E29').Value = '  -1.96%  '
$ws.Range('D30').Value = '0.06146'
$ws.Range('E30').Value = '  -7.77%  '
$ws.Range('D31').Value = '1.421'
$ws.Range('E31').Value = '  -1.30%  '
$ws.Range('D32').Value = '3.792'
$ws.Range('E32').Value = '  -0.92%  '
$ws.Range('D33').Value = '3.764'
$ws.Range('D34').Value = '1.728'
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('D35').Value = '1.053'
$ws.Range('E35').Value = '  -3.18%  '
$ws.Range('D36').Value = '0.6423'
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('D37').Value = '2.502'
$ws.Range('E37').Value = '  -1.82%  '
$ws.Range('D38').Value = '2.711'
$ws.Range('E38').Value = '  -0.86%  '
$ws.Range('D39').Value = '0.01691'
$ws.Range('E39').Value = '  -2.40%  '
$ws.Range('D40').Value = '1.147.52'
$ws.Range('E40').Value = '  -3.17%  '
$ws.Range('D41').Value = '6.324'
$ws.Range('E41').Value = '  -2.63%  '
$ws.Range('E42').Value = '  -3.08%  '
$ws.Range('E43').Value = '  -0.26%  '
$ws.Range('D44').Value = '100.60'
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('D45').Value = '1.947.15'
$ws.Range('E45').Value = '  -1.63%  '
$ws.Range('D46').Value = '60.22'
$ws.Range('E46').Value = '  -2.80%  '
$ws.Range('D47').Value = '0.00000000110'
$ws.Range('E47').Value = '  -4.51%  '
$ws.Range('D48').Value = '1.587'
$ws.Range('E48').Value = '  +1.26%  '
$ws.Range('D49').Value = '8.332'
$ws.Range('E49').Value = '  -1.67%  '
$ws.Range('D50').Value = '0.05452'
$ws.Range('E50').Value = '  -0.97%  '
$ws.Range('E51').Value = '  -1.86%  '
